# Updates crypto price/volume figures on sheet1 to reflect the latest
# GitHub Actions scrape, as captured in the commit
# "Updated cryptos list on Thu Apr 20 17:54:33 UTC 2023 with GitHub Actions".
#
# Price (column D) and Volume(1h) (column E) cells store plain text in the
# workbook (many prices use locale-style thousands separators such as
# "28.818.04"). For the handful of prices that also happen to look like
# ordinary decimal numbers to Excel (e.g. "323.24"), we briefly force the
# cell to text format before assigning the value so Excel keeps it as a
# string instead of silently re-parsing it as a floating point number,
# then restore the default "Normal" style so no stray formatting is left
# behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.818.04"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "1.961.87"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4770"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.84%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4034"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08475"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.060"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.31%  "
$ws.Range("D13").Value = "1.983.62"
$ws.Range("E13").Value = "  -1.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.640"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.230"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.015"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.15%  "
$ws.Range("E18").Value = "  -3.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06593"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.793"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "28.840.32"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.41%  "
$ws.Range("D26").Value = "2.209.45"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.948"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.153"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09590"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.459"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.677"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.678"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02358"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.268"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06208"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.742"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6222"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1918"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.341"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5956"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.072"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.427"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000335"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06834"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.21%  "
